# This edit rotates the observation data contained in rows 2, 3, 4 and 6
# of the "Artfynd" sheet: the content that used to live in row 2 moves to
# row 6, row 6's content moves to row 3, row 3's content moves to row 4,
# and row 4's content moves to row 2. All other columns (location/
# observer metadata) are identical across these four rows, so only the
# species-specific columns (A, B, E, F, G, H, I, J, Q, R) actually change
# value, plus the (empty) "Kön" column L which is present for rows 3/4
# and absent for rows 2/6 before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 4 (Murgröna / Hedera helix)
$ws.Range("A2").Value = 111412870
$ws.Range("B2").Value = 108022
$ws.Range("E2").Value = 219677
$ws.Range("F2").Value = "Murgröna"
$ws.Range("G2").Value = "Hedera helix"
$ws.Range("H2").Value = "L."
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("Q2").Value = 493001.1390786725
$ws.Range("R2").Value = 6227751.92766118

# Row 3 <- old row 6 (Liten stinksvamp / Mutinus caninus)
$ws.Range("A3").Value = 111413026
$ws.Range("B3").Value = 89007
$ws.Range("E3").Value = 1068
$ws.Range("F3").Value = "Liten stinksvamp"
$ws.Range("G3").Value = "Mutinus caninus"
$ws.Range("H3").Value = "(Schaeff.:Pers.) Fr."
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("L3").Value = ""
$ws.Range("Q3").Value = 492828.1855948549
$ws.Range("R3").Value = 6227996.970613244

# Row 4 <- old row 3 (Myskmadra / Galium odoratum)
$ws.Range("A4").Value = 111412858
$ws.Range("B4").Value = 103369
$ws.Range("E4").Value = 221423
$ws.Range("F4").Value = "Myskmadra"
$ws.Range("G4").Value = "Galium odoratum"
$ws.Range("H4").Value = "(L.) Scop."
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("Q4").Value = 493001.1390786725
$ws.Range("R4").Value = 6227751.92766118

# Row 6 <- old row 2 (Fjällsopp / Strobilomyces strobilaceus)
$ws.Range("A6").Value = 111412931
$ws.Range("B6").Value = 88680
$ws.Range("E6").Value = 1541
$ws.Range("F6").Value = "Fjällsopp"
$ws.Range("G6").Value = "Strobilomyces strobilaceus"
$ws.Range("H6").Value = "(Scop.:Fr.) Berk."
$ws.Range("I6").Value = "5"
$ws.Range("J6").Value = "fruktkroppar"
$ws.Range("Q6").Value = 492929.9770629382
$ws.Range("R6").Value = 6227867.193217421
